$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.552.69'
$ws.Range("D3").Value = '2.003.06'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("E4").Value = '  +1.19%  '
$ws.Range("D5").Value = '''329.44'
$ws.Range("E5").Value = '  -4.01%  '
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("D7").Value = '''0.5000'
$ws.Range("E7").Value = '  -4.53%  '
$ws.Range("D8").Value = '''0.4217'
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("D9").Value = '''54.35'
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").Value = '''0.09009'
$ws.Range("E10").Value = '  -3.42%  '
$ws.Range("E11").Value = '  -4.40%  '
$ws.Range("E12").Value = '  -6.47%  '
$ws.Range("D13").Value = '2.057.28'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '''8.026'
$ws.Range("E14").Value = '  -6.68%  '
$ws.Range("D15").Value = '''6.465'
$ws.Range("E15").Value = '  -6.26%  '
$ws.Range("D16").Value = '''1.014'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").Value = '''94.39'
$ws.Range("E17").Value = '  -6.78%  '
$ws.Range("E18").Value = '  -3.96%  '
$ws.Range("D19").Value = '''0.06676'
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").Value = '''19.65'
$ws.Range("E20").Value = '  -7.01%  '
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("D22").Value = '''5.967'
$ws.Range("E22").Value = '  -5.83%  '
$ws.Range("D23").Value = '29.598.75'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("E24").Value = '  -4.46%  '
$ws.Range("D25").Value = '''2.302'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '''159.01'
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("D27").Value = '''20.72'
$ws.Range("E27").Value = '  -4.98%  '
$ws.Range("D28").Value = '''6.347'
$ws.Range("E28").Value = '  -5.38%  '
$ws.Range("E29").Value = '  -8.83%  '
$ws.Range("D30").Value = '''128.18'
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("E31").Value = '  -7.21%  '
$ws.Range("D32").Value = '''0.09958'
$ws.Range("E32").Value = '  -4.75%  '
$ws.Range("D33").Value = '''1.565'
$ws.Range("E33").Value = '  -6.17%  '
$ws.Range("D34").Value = '''5.828'
$ws.Range("E34").Value = '  -6.59%  '
$ws.Range("D35").Value = '''3.802'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  -6.22%  '
$ws.Range("D37").Value = '''9.269'
$ws.Range("E37").Value = '  -8.84%  '
$ws.Range("D38").Value = '''0.06424'
$ws.Range("E38").Value = '  -6.16%  '
$ws.Range("D39").Value = '''1.306'
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("D40").Value = '''0.6537'
$ws.Range("E40").Value = '  -6.43%  '
$ws.Range("E41").Value = '  -6.83%  '
$ws.Range("D42").Value = '''0.2046'
$ws.Range("E42").Value = '  -7.58%  '
$ws.Range("D43").Value = '''1.012'
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").Value = '''0.6350'
$ws.Range("E44").Value = '  -6.96%  '
$ws.Range("D45").Value = '''13.42'
$ws.Range("E45").Value = '  -6.50%  '
$ws.Range("E46").Value = '  -6.14%  '
$ws.Range("D47").Value = '''1.304'
$ws.Range("E47").Value = '  -5.11%  '
$ws.Range("D48").Value = '''3.513'
$ws.Range("E48").Value = '  -3.34%  '
$ws.Range("D49").Value = '''0.00000000334'
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("D50").Value = '''0.06990'
$ws.Range("E50").Value = '  -3.28%  '
$ws.Range("D51").Value = '''1.128'
$ws.Range("E51").Value = '  -6.78%  '
